# Apply updated crypto market data (prices and 1h volume changes) as
# scraped by the GitHub Actions job on Mon Aug 28 13:28:22 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column holds plain text (not numbers), so keep it formatted
# as text before writing values that could otherwise be auto-coerced into
# numeric cells by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.230.94'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '1.659.18'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  -0.64%  '
$ws.Range('D5').Value = '218.59'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('D6').Value = '0.5237'
$ws.Range('E6').Value = '  -2.09%  '
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('D9').Value = '0.06312'
$ws.Range('E9').Value = '  -1.20%  '
$ws.Range('D10').Value = '20.63'
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('D11').Value = '0.07789'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('D13').Value = '1.650.36'
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('D14').Value = '1.887.95'
$ws.Range('E14').Value = '  -0.45%  '
$ws.Range('D15').Value = '0.5628'
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').Value = '0.0₅8057'
$ws.Range('E16').Value = '  -1.57%  '
$ws.Range('D17').Value = '65.17'
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D18').Value = '26.222.71'
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('E20').Value = '  +0.50%  '
$ws.Range('D21').Value = '194.30'
$ws.Range('E21').Value = '  -0.43%  '
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('D23').Value = '6.019'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('D24').Value = '1.005'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').Value = '145.41'
$ws.Range('E25').Value = '  -0.71%  '
$ws.Range('D26').Value = '0.1208'
$ws.Range('E26').Value = '  -1.50%  '
$ws.Range('D27').Value = '7.217'
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('D28').Value = '16.02'
$ws.Range('E28').Value = '  -0.91%  '
$ws.Range('D29').Value = '1.498'
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('D30').Value = '0.05649'
$ws.Range('E30').Value = '  -3.49%  '
$ws.Range('E31').Value = '  -0.83%  '
$ws.Range('D32').Value = '3.476'
$ws.Range('E32').Value = '  -3.05%  '
$ws.Range('D33').Value = '3.360'
$ws.Range('E33').Value = '  +2.08%  '
$ws.Range('D34').Value = '1.602'
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('D36').Value = '2.404'
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('D37').Value = '0.9428'
$ws.Range('D38').Value = '0.5752'
$ws.Range('E38').Value = '  -1.29%  '
$ws.Range('D39').Value = '0.01604'
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('D40').Value = '5.995'
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('B41').Value = 'mCoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D41').Value = '2.574'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.051.88'
$ws.Range('E42').Value = '  -2.25%  '
$ws.Range('D43').Value = '0.8463'
$ws.Range('E43').Value = '  -2.14%  '
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('D45').Value = '102.89'
$ws.Range('E45').Value = '  -1.22%  '
$ws.Range('D46').Value = '1.799.04'
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('D48').Value = '0.0₈105'
$ws.Range('E48').Value = '  +1.67%  '
$ws.Range('E49').Value = '  -1.05%  '
$ws.Range('D50').Value = '0.05321'
$ws.Range('E50').Value = '  +2.96%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.4352'
$ws.Range('E51').Value = '  -0.89%  '
